$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C36").ClearContents()
$ws.Range("C37").ClearContents()
